# Update crypto price/volume data per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.371.36"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.882.25"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7124"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08012"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3164"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.12"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08331"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "1.900.48"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.274"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.93"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7188"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.372"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008671"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  +5.27%  "
$ws.Range("D18").Value = "29.397.00"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.06"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "2.154.90"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.853"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.099"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.37"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.447"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.355"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  -6.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05401"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.950"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7754"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.187"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01893"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "1.274.16"
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.749"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.525"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9192"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.17"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.58"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("D47").Value = "2.042.53"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.817"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5223"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.585"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4385"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +1.29%  "
